$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) The six transaction rows (2-7) used to all point at the "Conta
#    Inicial" account label; that label is retired and the same rows now
#    belong to "Conta 1".
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "Conta 1"
$ws.Range("F3").Value = "Conta 1"
$ws.Range("F4").Value = "Conta 1"
$ws.Range("F5").Value = "Conta 1"
$ws.Range("F6").Value = "Conta 1"
$ws.Range("F7").Value = "Conta 1"

# ---------------------------------------------------------------------
# 2) Duplicate the same six transactions for "Conta 2" (rows 8-13) and
#    "Conta 3" (rows 14-19). Copy formatting from the source block first
#    (so styles/number-formats line up exactly), then stamp in the
#    literal values.
# ---------------------------------------------------------------------
$src = $ws.Range("A2:F7")

$src.Copy()
$ws.Range("A8:F13").PasteSpecial(-4122) # xlPasteFormats

$src.Copy()
$ws.Range("A14:F19").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# Conta 2 block (rows 8-13 mirror rows 2-7)
$ws.Range("A8").Value = 43328
$ws.Range("B8").Value = "Salário"
$ws.Range("C8").Value = "A definir"
$ws.Range("D8").Value = 1200
$ws.Range("F8").Value = "Conta 2"

$ws.Range("A9").Value = 43328
$ws.Range("B9").Value = "Energia Elétrica"
$ws.Range("C9").Value = "Energia Elétrica"
$ws.Range("D9").Value = -189.4
$ws.Range("F9").Value = "Conta 2"

$ws.Range("A10").Value = 43328
$ws.Range("B10").Value = "Conta de Água"
$ws.Range("C10").Value = "Água"
$ws.Range("D10").Value = -94.5
$ws.Range("E10").Value = "conteúdo 1"
$ws.Range("F10").Value = "Conta 2"

$ws.Range("A11").Value = 43329
$ws.Range("B11").Value = "Aluguel"
$ws.Range("C11").Value = "Aluguel"
$ws.Range("D11").Value = -600
$ws.Range("F11").Value = "Conta 2"

$ws.Range("A12").Value = 43138
$ws.Range("B12").Value = "Condomínio"
$ws.Range("C12").Value = "Condomínio"
$ws.Range("D12").Value = -300
$ws.Range("E12").Value = "conteúdo 3"
$ws.Range("F12").Value = "Conta 2"

$ws.Range("A13").Value = 43138
$ws.Range("B13").Value = "Compras supermercado"
$ws.Range("C13").Value = "Alimentação"
$ws.Range("D13").Value = -384.5
$ws.Range("E13").Value = "conteúdo 2"
$ws.Range("F13").Value = "Conta 2"

# Conta 3 block (rows 14-19 mirror rows 2-7)
$ws.Range("A14").Value = 43328
$ws.Range("B14").Value = "Salário"
$ws.Range("C14").Value = "A definir"
$ws.Range("D14").Value = 1200
$ws.Range("F14").Value = "Conta 3"

$ws.Range("A15").Value = 43328
$ws.Range("B15").Value = "Energia Elétrica"
$ws.Range("C15").Value = "Energia Elétrica"
$ws.Range("D15").Value = -189.4
$ws.Range("F15").Value = "Conta 3"

$ws.Range("A16").Value = 43328
$ws.Range("B16").Value = "Conta de Água"
$ws.Range("C16").Value = "Água"
$ws.Range("D16").Value = -94.5
$ws.Range("E16").Value = "conteúdo 1"
$ws.Range("F16").Value = "Conta 3"

$ws.Range("A17").Value = 43329
$ws.Range("B17").Value = "Aluguel"
$ws.Range("C17").Value = "Aluguel"
$ws.Range("D17").Value = -600
$ws.Range("F17").Value = "Conta 3"

$ws.Range("A18").Value = 43138
$ws.Range("B18").Value = "Condomínio"
$ws.Range("C18").Value = "Condomínio"
$ws.Range("D18").Value = -300
$ws.Range("E18").Value = "conteúdo 3"
$ws.Range("F18").Value = "Conta 3"

$ws.Range("A19").Value = 43138
$ws.Range("B19").Value = "Compras supermercado"
$ws.Range("C19").Value = "Alimentação"
$ws.Range("D19").Value = -384.5
$ws.Range("E19").Value = "conteúdo 2"
$ws.Range("F19").Value = "Conta 3"

# ---------------------------------------------------------------------
# 3) Normalize all row heights (15.95 / 15.6 -> 15) for rows 1-19, and
#    extend the sheet with twelve more blank rows (20-31) at height 15.
# ---------------------------------------------------------------------
$ws.Range("A1:A31").RowHeight = 15
